# Reproduce the edits described in the commit diff for
# "CNN+att results/cnn+att_result_full.xlsx":
#   - row 12 summary formulas: LEFT(...,5) -> ROUND(...,3)
#   - remove the last column (AJ), which held a redundant run-timestamp
#     label per row (the sheet already has an identical "loss_tr"/"nan"
#     pair in columns AH/AI); deleting AJ also drops the now-orphaned
#     shared-string timestamps automatically
#   - minor view-state tweaks: selection moved to column AA, and columns
#     Y/AA get an explicit (default) width marker

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the row-12 "mean ± std" summary formulas -------------------
# Column B (the master of the original shared-formula group B12:AH12) is
# rewritten first so the whole shared group picks up the new formula text;
# column A is a standalone (non-shared) formula and is set separately so it
# does not get folded back into the B:AH shared group.
$ws.Range("B12:AH12").Formula = '=ROUND(AVERAGE(B2:B11),3) &"±"& ROUND(_xlfn.STDEV.P(B2:B11),3)'
$ws.Range("A12").Formula = '=ROUND(AVERAGE(A2:A11),3) &"±"& ROUND(_xlfn.STDEV.P(A2:A11),3)'

# --- 2. Delete column AJ ----------------------------------------------------
# AJ duplicated the "loss_tr" timestamp label already present in AI; removing
# it shifts nothing else and the sheet's shared-string table is compacted
# automatically (the orphaned "conv_attn_..." timestamp strings drop out,
# leaving just the "nan" label used by AI).
$ws.Columns("AJ").Delete()

# --- 3. View-state touch-ups -------------------------------------------------
# Columns Y and AA get touched (without materially changing their width) and
# the active selection moves to the top of column AA.
$ws.Columns("Y").ColumnWidth = 9.140625
$ws.Columns("AA").ColumnWidth = 9.140625

$null = $ws.Range("AA1:AA1048576").Select()
$excel.ActiveWindow.ScrollColumn = 9
$excel.ActiveWindow.ScrollRow = 1
